$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-234 all hold the serial date value 45190
# (2023-09-21) and need to be updated to 45192 (2023-09-23).
for ($row = 2; $row -le 234; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
